$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "snaks" -> "snacks"
$ws.Range("C2").Value = "snacks"

# Row 5 (id 4): update date, amount (now numeric), and description
$ws.Range("B5").Value = 241008
$ws.Range("D5").Value = 1200
$ws.Range("E5").Value = "prime"

# Row 6 (id 5): update date, tag, amount, and description
$ws.Range("B6").Value = 241031
$ws.Range("C6").Value = "accessories"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "200"
$ws.Range("E6").Value = "pen"
